# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H,I,J,K,L,M,N) for the affected leve rows across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting sheets, mirroring the upstream
# Universalis price pull. Blank-string assignments clear cells that no
# longer carry a value (e.g. a profit column that reverts to "not computed").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2500
$ws.Range("J64").Value = 2500
$ws.Range("L64").Value = 2500
$ws.Range("N64").Value = -2996

$ws.Range("H67").Value = 2500
$ws.Range("J67").Value = 2500
$ws.Range("L67").Value = 2500
$ws.Range("N67").Value = -4216

$ws.Range("H88").Value = 4999
$ws.Range("J88").Value = 4999
$ws.Range("L88").Value = 4999
$ws.Range("N88").Value = -5811

$ws.Range("H91").Value = 4999
$ws.Range("J91").Value = 4999
$ws.Range("L91").Value = 4999
$ws.Range("N91").Value = -7807

$ws.Range("H127").Value = 400
$ws.Range("I127").Value = 400
$ws.Range("K127").Value = 1200
$ws.Range("M127").Value = 3760

$ws.Range("H129").Value = 408.5
$ws.Range("J129").Value = 417
$ws.Range("L129").Value = 1251
$ws.Range("N129").Value = -11251

$ws.Range("H132").Value = 1213.2222
$ws.Range("I132").Value = 1159.8572
$ws.Range("K132").Value = 3479.5716
$ws.Range("M132").Value = -949.5715999999998

$ws.Range("H135").Value = 491.75
$ws.Range("I135").Value = 491.75
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4425.75
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = ""
$ws.Range("N135").Value = -1890.75

$ws.Range("H138").Value = 5852.6665
$ws.Range("I138").Value = 7121
$ws.Range("K138").Value = 21363
$ws.Range("M138").Value = -16223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = ""
$ws.Range("N58").Value = 0

$ws.Range("H61").Value = 2025.3334
$ws.Range("I61").Value = 2348.4
$ws.Range("J61").Value = 410
$ws.Range("K61").Value = 2348.4
$ws.Range("L61").Value = 410
$ws.Range("M61").Value = -2136.4
$ws.Range("N61").Value = -834

$ws.Range("H74").Value = 2054.4285
$ws.Range("I74").Value = 2074.9
$ws.Range("K74").Value = 2074.9
$ws.Range("M74").Value = -1200.9

$ws.Range("H77").Value = 2054.4285
$ws.Range("I77").Value = 2074.9
$ws.Range("K77").Value = 10374.5
$ws.Range("M77").Value = -6006.5

$ws.Range("H132").Value = 92660.09
$ws.Range("I132").Value = 112655.11
$ws.Range("K132").Value = 337965.33
$ws.Range("M132").Value = -335435.33

$ws.Range("H136").Value = 2025.3334
$ws.Range("I136").Value = 2348.4
$ws.Range("J136").Value = 410
$ws.Range("K136").Value = 7045.200000000001
$ws.Range("L136").Value = 1230
$ws.Range("M136").Value = -4495.200000000001
$ws.Range("N136").Value = -6330

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 744.625
$ws.Range("I80").Value = 913.75
$ws.Range("J80").Value = 575.5
$ws.Range("K80").Value = 913.75
$ws.Range("L80").Value = 575.5
$ws.Range("M80").Value = 84.25
$ws.Range("N80").Value = -2571.5

$ws.Range("H83").Value = 744.625
$ws.Range("I83").Value = 913.75
$ws.Range("J83").Value = 575.5
$ws.Range("K83").Value = 4568.75
$ws.Range("L83").Value = 2877.5
$ws.Range("M83").Value = 423.25
$ws.Range("N83").Value = -12861.5

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 15625
$ws.Range("I60").Value = 15625
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 15625
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -15114

$ws.Range("H74").Value = 45332.5
$ws.Range("J74").Value = 45332.5
$ws.Range("L74").Value = 45332.5
$ws.Range("N74").Value = -47080.5

$ws.Range("H77").Value = 45332.5
$ws.Range("J77").Value = 45332.5
$ws.Range("L77").Value = 135997.5
$ws.Range("N77").Value = -144733.5

$ws.Range("H132").Value = 3459.524
$ws.Range("I132").Value = 2797.3684
$ws.Range("K132").Value = 8392.1052
$ws.Range("M132").Value = -5862.1052

$ws.Range("H134").Value = 770.55554
$ws.Range("I134").Value = 641.875
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 1925.625
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = 609.375
$ws.Range("N134").Value = -10470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3725.7144
$ws.Range("I5").Value = 3681.3333
$ws.Range("J5").Value = 3992
$ws.Range("K5").Value = 11043.9999
$ws.Range("L5").Value = 11976
$ws.Range("M5").Value = -10931.9999
$ws.Range("N5").Value = -12200

$ws.Range("H55").Value = 4000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = 12000
$ws.Range("N55").Value = -12354

$ws.Range("H92").Value = 521.375
$ws.Range("I92").Value = 485.14285
$ws.Range("J92").Value = 775
$ws.Range("K92").Value = 1455.42855
$ws.Range("L92").Value = 2325
$ws.Range("M92").Value = -207.4285500000001
$ws.Range("N92").Value = -4821

$ws.Range("H135").Value = 3725.7144
$ws.Range("I135").Value = 3681.3333
$ws.Range("J135").Value = 3992
$ws.Range("K135").Value = 33131.9997
$ws.Range("L135").Value = 35928
$ws.Range("M135").Value = -30596.9997
$ws.Range("N135").Value = -40998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2399.182
$ws.Range("I132").Value = 2173.875
$ws.Range("K132").Value = 6521.625
$ws.Range("M132").Value = -3991.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 13749.375
$ws.Range("I35").Value = 7250
$ws.Range("J35").Value = 20248.75
$ws.Range("K35").Value = 7250
$ws.Range("L35").Value = 20248.75
$ws.Range("M35").Value = -6914
$ws.Range("N35").Value = -20920.75

$ws.Range("H132").Value = 557711.25
$ws.Range("I132").Value = 627049.5
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 1881148.5
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -1878618.5
$ws.Range("N132").Value = -14075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""

$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

$ws.Range("H132").Value = 1057.2858
$ws.Range("I132").Value = 907.46155
$ws.Range("K132").Value = 2722.38465
$ws.Range("M132").Value = -192.38465
